# Add calculation of total return
# Populates Capital Gains (H), Dividends Paid (I) and Total Return (J) for
# the existing portfolio rows, and refreshes the Market Price / Annual
# Dividend per Share inputs (and a couple of other values for the SCHD
# position) that feed those calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Atlantica Sustainable Infrastructure (AY), 687 shares ---
$ws.Range("E2").Value = 21.6
$ws.Range("F2").Value = 1.78
$ws.Range("H2").Value = -8.279999999999999
$ws.Range("I2").Value = 611.4300000000001
$ws.Range("J2").Value = -4.50106157112526

# --- Row 3: Atlantica Sustainable Infrastructure (AY), 400 shares ---
$ws.Range("E3").Value = 21.6
$ws.Range("F3").Value = 1.78
$ws.Range("H3").Value = -12.41
$ws.Range("I3").Value = 356
$ws.Range("J3").Value = -8.799675587996749

# --- Row 4: Schwab US Dividend Equity ETF (SCHD) ---
$ws.Range("C4").Value = 983
$ws.Range("D4").Value = 75.65000000000001
$ws.Range("E4").Value = 75.81999999999999
$ws.Range("F4").Value = 2.66
$ws.Range("G4").Value = "27.07.2023"
$ws.Range("H4").Value = 0.22
$ws.Range("I4").Value = 1376.2
$ws.Range("J4").Value = 2.075346992729667
